$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the existing "DAM 75" block header (row 102) to "TAWATAIA DAM"
# ---------------------------------------------------------------------------
$ws.Range("D102").Value = "TAWATAIA DAM"

# ---------------------------------------------------------------------------
# 2. Append a brand-new "DAM 75" coordinate block starting at row 124
# ---------------------------------------------------------------------------
$headerRow = 124

$ws.Cells.Item($headerRow, 1).Value = "X"
$ws.Cells.Item($headerRow, 2).Value = "Y"
$ws.Cells.Item($headerRow, 3).Value = "Z"
$ws.Cells.Item($headerRow, 4).Value = "DAM 75"
$ws.Cells.Item($headerRow, 5).Value = "gap between points"
$ws.Cells.Item($headerRow, 6).Value = "running distance"
$ws.Cells.Item($headerRow, 7).Value = "RADAR LEVEL : "
$ws.Cells.Item($headerRow, 8).Value = "PAVER LEVEL"
$ws.Cells.Item($headerRow, 9).Value = "CULVERT INVERT"

# match formatting of the header cells to the other dam-block headers
$ws.Cells.Item($headerRow, 4).Style = $ws.Cells.Item(102, 4).Style
$ws.Cells.Item($headerRow, 5).Style = $ws.Cells.Item(102, 5).Style
$ws.Cells.Item($headerRow, 6).Style = $ws.Cells.Item(102, 6).Style

# Row data: row, X, Y, Z
$rows = @(
    @(125, 1837443.1, 5496185.5999999996, 171.98400000000001),
    @(126, 1837447.15, 5496194.8399999999, 171.976),
    @(127, 1837451.26, 5496204.3499999996, 171.989),
    @(128, 1837455.19, 5496213.7999999998, 172.02600000000001),
    @(129, 1837459.09, 5496223.3099999996, 172.04),
    @(130, 1837463.07, 5496233.0599999996, 172.078),
    @(131, 1837466.97, 5496243.2699999996, 172.119),
    @(132, 1837470.37, 5496253.1900000004, 172.14500000000001),
    @(133, 1837473.48, 5496263.2000000002, 172.148),
    @(134, 1837476.06, 5496273.5, 172.18299999999999),
    @(135, 1837478.42, 5496283.8899999997, 172.21700000000001),
    @(136, 1837480.01, 5496294.7699999996, 172.203),
    @(137, 1837481.36, 5496305.5700000003, 172.17400000000001),
    @(138, 1837482.11, 5496315.9500000002, 172.12799999999999),
    @(139, 1837483.51, 5496328.2300000004, 172.17),
    @(140, 1837484.22, 5496338.5, 172.33500000000001),
    @(141, 1837484.94, 5496348.2300000004, 172.64),
    @(142, 1837485.49, 5496358.3799999999, 173.21799999999999),
    @(143, 1837486.58, 5496368.5800000001, 174.018),
    @(144, 1837487.59, 5496378.7000000002, 174.75800000000001),
    @(145, 1837489.15, 5496388.9100000001, 175.10599999999999),
    @(146, 1837490.7, 5496399.6699999999, 174.86600000000001),
    @(147, 1837491.89, 5496409.7000000002, 174.44499999999999),
    @(148, 1837493.05, 5496419.4699999997, 174.096),
    @(149, 1837493.96, 5496429.5, 173.79300000000001),
    @(150, 1837494.77, 5496439.6200000001, 173.52),
    @(151, 1837495.48, 5496449.9299999997, 173.31800000000001),
    @(152, 1837496.2, 5496459.96, 173.16200000000001),
    @(153, 1837496.83, 5496469.9900000002, 172.958),
    @(154, 1837497.45, 5496479.96, 172.77600000000001),
    @(155, 1837498.13, 5496490.1100000003, 172.57900000000001),
    @(156, 1837498.86, 5496500.4500000002, 172.441),
    @(157, 1837499.64, 5496510.7999999998, 172.37),
    @(158, 1837500.38, 5496521.29, 172.333),
    @(159, 1837501.21, 5496531.6100000003, 172.37100000000001),
    @(160, 1837501.9, 5496541.6299999999, 172.43299999999999),
    @(161, 1837502.51, 5496552.0800000001, 172.54300000000001)
)

$firstDataRow = $rows[0][0]
$lastDataRow = $rows[$rows.Length - 1][0]

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = [double]$r[1]
    $ws.Cells.Item($rowNum, 2).Value = [double]$r[2]
    $ws.Cells.Item($rowNum, 3).Value = [double]$r[3]
}

# First data row: start-of-run distance + radar/paver/culvert levels
$ws.Cells.Item($firstDataRow, 6).Value = 0
$ws.Cells.Item($firstDataRow, 9).Value = 168.98
$ws.Cells.Item($firstDataRow, 8).Formula = "=I$firstDataRow+4"
$ws.Cells.Item($firstDataRow, 7).Formula = "=I$firstDataRow+8"

# Remaining rows: running SQRT distance + cumulative chainage
for ($rn = $firstDataRow + 1; $rn -le $lastDataRow; $rn++) {
    $prev = $rn - 1
    $ws.Cells.Item($rn, 5).Formula = "=SQRT((A$rn-A$prev)^2 + (B$prev-B$rn)^2)"
    $ws.Cells.Item($rn, 6).Formula = "=F$prev+E$rn"
}

Write-Host "done"
